# ResearchDiary/RelatedWorkTabelle.xlsx - add new related-work entries to the
# "General" sheet table (IPQ / embodiment questionnaire + perspective-switch
# papers), wire up their source-link hyperlinks, extend the table range and
# add hyperlinks to a handful of pre-existing rows that only had bare URL
# text in the "Comments" column so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# xlPasteFormats / xlPasteSpecial constants used below
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Add hyperlinks to existing rows whose "Comments" (J) cell already holds
#    the raw URL as text but never got turned into an actual hyperlink.
#    We reuse the exact look of J2 (an existing, properly-styled hyperlink
#    cell) via copy/paste-special so the cell keeps the normal "Hyperlink"
#    look instead of picking up a fresh ad-hoc style.
# ---------------------------------------------------------------------------
$hyperlinkStyleSource = $ws.Cells.Item(2, 10)

$rowsNeedingLinks = @(42, 58, 59, 57, 31, 4)
foreach ($r in $rowsNeedingLinks) {
    $cell = $ws.Cells.Item($r, 10)
    $url = $cell.Text
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    $hyperlinkStyleSource.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# 2) Append the four new rows (60-63) to the table.
# ---------------------------------------------------------------------------

# Row 60 - quick IPQ questionnaire reference (only Paper/Summary/Comments set)
$ws.Cells.Item(60, 2).Value = "ipq"
$ws.Cells.Item(60, 9).Value = "presence questionnaire used by hoppe for perspective continuum"
$c = $ws.Cells.Item(60, 10)
$c.Value = "https://www.igroup.org/pq/ipq/download.php"
$ws.Hyperlinks.Add($c, "https://www.igroup.org/pq/ipq/download.php") | Out-Null
$hyperlinkStyleSource.Copy()
$c.PasteSpecial($xlPasteFormats)

# Row 61 - Gorisse et al. (embodiment questionnaire)
$ws.Cells.Item(61, 1).Value = "Geoffrey gorisse et al."
$titleCell = $ws.Cells.Item(61, 2)
$titleCell.Value = "First- and Third-Person Perspectives in Immersive Virtual Environments: Presence and Performance Analysis of Embodied Users"
$titleCell.Font.Size = 24
$titleCell.VerticalAlignment = -4108
$titleCell.WrapText = $false
$ws.Cells.Item(61, 4).Value = 2017
$ws.Cells.Item(61, 9).Value = "embodiment questionnaire used by hoppe for perspective continuum"
$c = $ws.Cells.Item(61, 10)
$c.Value = "https://www.frontiersin.org/journals/robotics-and-ai/articles/10.3389/frobt.2017.00033/full"
$ws.Hyperlinks.Add($c, "https://www.frontiersin.org/journals/robotics-and-ai/articles/10.3389/frobt.2017.00033/full") | Out-Null
$hyperlinkStyleSource.Copy()
$c.PasteSpecial($xlPasteFormats)
$ws.Rows.Item(61).RowHeight = 31

# Row 62 - ShadowClones paper (title style copied from row 57's title cell)
$ws.Cells.Item(62, 1).Value = "kazuma takada et al."
$titleCell62 = $ws.Cells.Item(62, 2)
$titleCell62.Value = "ShadowClones: an Interface to Maintain a Multiple Sense of Body-space Coordination in Multiple Visual Perspectives"
$ws.Cells.Item(57, 2).Copy()
$titleCell62.PasteSpecial($xlPasteFormats)
$ws.Cells.Item(62, 4).Value = 2023
$ws.Cells.Item(62, 9).Value = "shadow clones, user controls multiple things at the same time"
$ws.Cells.Item(62, 10).Value = "https://dl.acm.org/doi/10.1145/3582700.3582706"
$ws.Rows.Item(62).RowHeight = 180

# Row 63 - "I'm Transforming!" paper (title style copied from row 58's title cell)
$ws.Cells.Item(63, 1).Value = "riku otono et al."
$titleCell63 = $ws.Cells.Item(63, 2)
$titleCell63.Value = "I'm Transforming! Effects of Visual Transitions to Change of Avatar on the Sense of Embodiment in AR"
$ws.Cells.Item(58, 2).Copy()
$titleCell63.PasteSpecial($xlPasteFormats)
$ws.Cells.Item(63, 4).Value = 2023
$ws.Cells.Item(63, 9).Value = "paper on switching embodiment and how to do it smoothly, they link it to physical action"
$ws.Cells.Item(63, 10).Value = "https://ieeexplore.ieee.org/document/10108469"
$ws.Rows.Item(63).RowHeight = 60

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Extend the table / autofilter range to cover the new rows.
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J63"))

# ---------------------------------------------------------------------------
# 4) View-state bits: selection moved to L62, scrolled so B1 leads, sheet
#    frozen pane still anchored at row 57 when last saved.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("L62").Select()
